$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("10 yrs Old_ConnectNearBy"): move the selection from B21 to A2 ---
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Range("A2").Select()

# --- Sheet 3 ("Sheet3"): fill in the new results table ---
$ws3 = $wb.Worksheets.Item(3)

# Header row (new shared strings get created the first time each text is used)
$ws3.Cells.Item(1, 1).Value = "leaves"
$ws3.Cells.Item(1, 2).Value = "Diff.tree"
$ws3.Cells.Item(1, 3).Value = "intern_br_sim"
$ws3.Cells.Item(1, 4).Value = "intern_br_true"
$ws3.Cells.Item(1, 5).Value = "diff.tree/2*intern_br_sim"
$ws3.Cells.Item(1, 6).Value = "diff.tree/2*intern_br_all"

# Data rows
$ws3.Cells.Item(2, 1).Value = 5
$ws3.Cells.Item(2, 2).Value = 4
$ws3.Cells.Item(2, 3).Value = 3
$ws3.Cells.Item(2, 4).Value = 3
$ws3.Cells.Item(2, 5).Value = 0.66
$ws3.Cells.Item(2, 6).Value = 0.33

$ws3.Cells.Item(3, 1).Value = 10
$ws3.Cells.Item(3, 2).Value = 16
$ws3.Cells.Item(3, 3).Value = 8
$ws3.Cells.Item(3, 4).Value = 8
$ws3.Cells.Item(3, 5).Value = 1
$ws3.Cells.Item(3, 6).Value = 0.5

$ws3.Cells.Item(4, 1).Value = 15
$ws3.Cells.Item(4, 2).Value = 26
$ws3.Cells.Item(4, 3).Value = 13
$ws3.Cells.Item(4, 4).Value = 13
$ws3.Cells.Item(4, 5).Value = 1
$ws3.Cells.Item(4, 6).Value = 0.5

$ws3.Cells.Item(5, 1).Value = 20
$ws3.Cells.Item(5, 2).Value = 36
$ws3.Cells.Item(5, 3).Value = 18
$ws3.Cells.Item(5, 4).Value = 18
$ws3.Cells.Item(5, 5).Value = 1

$ws3.Cells.Item(6, 1).Value = 25
$ws3.Cells.Item(6, 2).Value = 46
$ws3.Cells.Item(6, 3).Value = 23
$ws3.Cells.Item(6, 4).Value = 23
$ws3.Cells.Item(6, 5).Value = 1

$ws3.Cells.Item(7, 1).Value = 30
$ws3.Cells.Item(7, 2).Value = 56
$ws3.Cells.Item(7, 3).Value = 28
$ws3.Cells.Item(7, 4).Value = 28
$ws3.Cells.Item(7, 5).Value = 1

$ws3.Cells.Item(8, 1).Value = 35
$ws3.Cells.Item(8, 2).Value = 66
$ws3.Cells.Item(8, 3).Value = 33
$ws3.Cells.Item(8, 4).Value = 33
$ws3.Cells.Item(8, 5).Value = 1

$ws3.Cells.Item(9, 1).Value = 40
$ws3.Cells.Item(9, 2).Value = 76
$ws3.Cells.Item(9, 3).Value = 38
$ws3.Cells.Item(9, 4).Value = 38
$ws3.Cells.Item(9, 5).Value = 1

$ws3.Cells.Item(10, 1).Value = 45
$ws3.Cells.Item(10, 2).Value = 86
$ws3.Cells.Item(10, 3).Value = 43
$ws3.Cells.Item(10, 4).Value = 43
$ws3.Cells.Item(10, 5).Value = 1

$ws3.Cells.Item(11, 1).Value = 50
$ws3.Cells.Item(11, 2).Value = 92
$ws3.Cells.Item(11, 3).Value = 48
$ws3.Cells.Item(11, 4).Value = 48
$ws3.Cells.Item(11, 5).Value = 0.95

$ws3.Cells.Item(12, 1).Value = 55
$ws3.Cells.Item(12, 2).Value = 106
$ws3.Cells.Item(12, 3).Value = 53
$ws3.Cells.Item(12, 4).Value = 53
$ws3.Cells.Item(12, 5).Value = 1

$ws3.Cells.Item(13, 1).Value = 60
$ws3.Cells.Item(13, 2).Value = 116
$ws3.Cells.Item(13, 3).Value = 58
$ws3.Cells.Item(13, 4).Value = 58
$ws3.Cells.Item(13, 5).Value = 1

$ws3.Cells.Item(14, 1).Value = 65
$ws3.Cells.Item(14, 2).Value = 126
$ws3.Cells.Item(14, 3).Value = 63
$ws3.Cells.Item(14, 4).Value = 63
$ws3.Cells.Item(14, 5).Value = 1

$ws3.Cells.Item(15, 1).Value = 70
$ws3.Cells.Item(15, 2).Value = 136
$ws3.Cells.Item(15, 3).Value = 68
$ws3.Cells.Item(15, 4).Value = 68
$ws3.Cells.Item(15, 5).Value = 1

$ws3.Cells.Item(16, 1).Value = 72
$ws3.Cells.Item(16, 2).Value = 140
$ws3.Cells.Item(16, 3).Value = 72
$ws3.Cells.Item(16, 4).Value = 72
$ws3.Cells.Item(16, 5).Value = 0.97
$ws3.Cells.Item(16, 6).Value = 0.48

# Column widths for the new columns E, F, G (values chosen so the lossy
# ColumnWidth -> pixel -> stored-width round trip lands as close as possible
# to the target stored widths of 25.8928571428571 / 26.015306122449 / 16.6683673469388)
$ws3.Columns.Item(5).ColumnWidth = 25.009523809523767
$ws3.Columns.Item(6).ColumnWidth = 25.131972789115665
$ws3.Columns.Item(7).ColumnWidth = 15.785034013605467

# Make Sheet3 the active sheet/tab, with F5 selected
[void]$ws3.Activate()
[void]$ws3.Range("F5").Select()
